$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: existing "2,4-D" row - B10 switches from a number to the "na" label ---
$ws.Range("B10").Value = "na"

# --- New pesticide rows (11-17) ---
# Column A values are entered in this particular order so that the shared-string
# table ends up with the same ordering as the authored workbook.
$ws.Range("A11").Value = "Sulfuric Acid"
$ws.Range("A13").Value = "Propanil"
$ws.Range("A12").Value = "Copper Sulfate Tribasic"
$ws.Range("A14").Value = "Chloropicrin"
$ws.Range("A15").Value = "Methyl Bromide (Bromomethane)"
$ws.Range("A16").Value = "Cyanamide"
$ws.Range("A17").Value = "Allyl Isothiocyanate"

$ws.Range("B11").Value = "na"
$ws.Range("C11").Value = 2140
$ws.Range("D11").Value = 34014

$ws.Range("B12").Value = "na"
$ws.Range("C12").Value = 2500
$ws.Range("D12").Value = 20584

$ws.Range("B13").Value = "na"
$ws.Range("C13").Value = 367
$ws.Range("D13").Value = 18083

$ws.Range("B14").Value = "na"
$ws.Range("C14").Value = 250
$ws.Range("D14").Value = 12923

$ws.Range("B15").Value = "na"
$ws.Range("C15").Value = 214
$ws.Range("D15").Value = 10742

$ws.Range("B16").Value = 388
$ws.Range("C16").Value = 125
$ws.Range("D16").Value = 9178

$ws.Range("B17").Value = 308
$ws.Range("C17").Value = 148
$ws.Range("D17").Value = 8908

# --- Column A is widened to fit the new, longer pesticide names ---
$ws.Columns.Item(1).ColumnWidth = 35.8

# --- Threaded comments added during review ---
$dComment = $ws.Range("D11").AddCommentThreaded("Could be incorrect")

$aComment = $ws.Range("A17").AddCommentThreaded("https://pubchem.ncbi.nlm.nih.gov/compound/5971#section=Toxicological-Information")
$aComment.AddReply("What to do if multiple reported LD50? Just average?")

# --- Leave the selection where Excel would land after entering the last row ---
$ws.Range("A18").Select()
